$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append (dates + numeric columns B:O)
$newRows = @(
    @("2020-12-24", 2132, 659, 987, 486, 2809, 366, 478, 1965, 30.91, 46.29, 22.8, 13.03, 17.02, 69.95),
    @("2020-12-25", 2133, 636, 927, 570, 2804, 352, 469, 1983, 29.82, 43.46, 26.72, 12.55, 16.73, 70.72),
    @("2020-12-26", 2134, 638, 908, 588, 2804, 360, 476, 1968, 29.9, 42.55, 27.55, 12.84, 16.98, 70.19),
    @("2020-12-27", 2136, 632, 950, 554, 2801, 364, 489, 1948, 29.59, 44.48, 25.94, 13, 17.46, 69.55),
    @("2020-12-28", 2137, 646, 962, 529, 2802, 363, 498, 1941, 30.23, 45.02, 24.75, 12.96, 17.77, 69.27)
)

$startRow = 303

# Scratch cell used to build the date text via a formula and then paste
# it back as a plain value. Going through Copy/PasteSpecial(values) avoids
# Excel's "smart" text-to-date autoconversion (and the style/number-format
# stamp that comes with it) that a direct .Value = "2020-12-24" assignment
# would trigger, keeping the cell as a plain shared-string cell like the
# rest of column A.
$scratch = $ws.Cells.Item(1000, 1)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]

    $scratch.Formula = '="' + $rowData[0] + '"'
    $scratch.Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4163)  # xlPasteValues

    for ($col = 2; $col -le 15; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }
}

$scratch.Clear()
$excel.CutCopyMode = 0
